# Birds sheet: fix bug with invalid input for bird id
# Adds the missing bird record (Bird ID 11) as row 12.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Birds")

$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "European Gouldian"
$ws.Range("C12").Value = "East Europe"
$ws.Range("D12").Value = "Male"
$ws.Range("E12").Value = 10

# Copy the date format from the row above so the new date cell reuses the
# existing "date" cell style instead of Excel minting a brand-new one.
$ws.Range("G11").Copy()
$ws.Range("G12").PasteSpecial(-4122)
$ws.Range("G12").Value = 45049

$ws.Range("H12").Value = "a3"
$ws.Range("I12").Value = 1
